# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 1054
    "F6"  = 3031
    "F7"  = 41
    "F8"  = 2121
    "F10" = 105
    "F11" = 970
    "F14" = 236
    "F17" = 36
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
